$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "TestCases": TC002 Run Mode -> No, add new TC003 row (Run Mode Yes,
# Results FAIL)
# ---------------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("TestCases")

# Row 3 (TC002) currently has style pattern A-C=s2, D=no-style; row 4 needs
# the exact same pattern, so copy it down first and then overwrite values.
# The old row 3 "Results" value (D3, FAIL) moves down to D4, so it must be
# cleared from its original location afterwards.
$wsCases.Range("A3:D3").Copy($wsCases.Range("A4"))
$wsCases.Range("D3").ClearContents()

$wsCases.Range("C3").Value() = "No"

$wsCases.Range("A4").Value() = "TC003"
$wsCases.Range("B4").Value() = "Add a product to cart"
$wsCases.Range("C4").Value() = "Yes"
$wsCases.Range("D4").Value() = "FAIL"

# ---------------------------------------------------------------------------
# Sheet "TestSteps": append the TC003 ("Add a product to cart") steps
# ---------------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("TestSteps")

# Rows 20-28 reuse the same visual style as the TC001 block (rows 2-10);
# copying that block over preserves per-cell styles (incl. the H column,
# which intentionally carries no explicit style).
$wsSteps.Range("A2:H10").Copy($wsSteps.Range("A20"))

# Relabel the copied block as TC003 steps.
$wsSteps.Range("A20:A28").Value() = "TC003"
$wsSteps.Range("B20").Value() = "TC003_01"
$wsSteps.Range("B21").Value() = "TC003_02"
$wsSteps.Range("B22").Value() = "TC003_03"
$wsSteps.Range("B23").Value() = "TC003_04"
$wsSteps.Range("B24").Value() = "TC003_05"
$wsSteps.Range("B25").Value() = "TC003_06"
$wsSteps.Range("B26").Value() = "TC003_07"
$wsSteps.Range("B27").Value() = "TC003_08"
$wsSteps.Range("B28").Value() = "TC003_09"

# Rows 27-28 (copied from TC001's "Click on My Account Link" / "Input
# User Name") actually describe the shopping-page steps, so overwrite them.
$wsSteps.Range("C27").Value() = "Click on Product Category"
$wsSteps.Range("D27").Value() = "ShoppingPage"
$wsSteps.Range("E27").Value() = "lnk_ProductCategory"
$wsSteps.Range("F27").Value() = "click"
$wsSteps.Range("G27").Value() = ""

$wsSteps.Range("C28").Value() = "Click on iPhones"
$wsSteps.Range("D28").Value() = "ShoppingPage"
$wsSteps.Range("E28").Value() = "lnk_iPhones"
$wsSteps.Range("F28").Value() = "click"
$wsSteps.Range("G28").Value() = ""
$wsSteps.Range("H28").Value() = "FAIL"

# Rows 29-33 are brand new steps with no H-column "Results" value; copy the
# A27:G27 style (s=4, no value) down as a formatting template for them.
$wsSteps.Range("A27:G27").Copy($wsSteps.Range("A29"))
$wsSteps.Range("A27:G27").Copy($wsSteps.Range("A30"))
$wsSteps.Range("A27:G27").Copy($wsSteps.Range("A31"))
$wsSteps.Range("A27:G27").Copy($wsSteps.Range("A32"))
$wsSteps.Range("A27:G27").Copy($wsSteps.Range("A33"))

$wsSteps.Range("A29:A33").Value() = "TC003"

$wsSteps.Range("B29").Value() = "TC003_10"
$wsSteps.Range("C29").Value() = "Verify and store product price"
$wsSteps.Range("D29").Value() = "ShoppingPage"
$wsSteps.Range("E29").Value() = "lbl_CurrentPrice"
$wsSteps.Range("F29").Value() = "storeValue"
$wsSteps.Range("G29").Value() = ""

$wsSteps.Range("B30").Value() = "TC003_11"
$wsSteps.Range("C30").Value() = "Add product to Cart"
$wsSteps.Range("D30").Value() = "ShoppingPage"
$wsSteps.Range("E30").Value() = "btn_AddToCart"
$wsSteps.Range("F30").Value() = "click"
$wsSteps.Range("G30").Value() = ""

$wsSteps.Range("B31").Value() = "TC003_12"
$wsSteps.Range("C31").Value() = "Navigate to Cart"
$wsSteps.Range("D31").Value() = "ShoppingPage"
$wsSteps.Range("E31").Value() = "btn_GoToCheckOut"
$wsSteps.Range("F31").Value() = "click"
$wsSteps.Range("G31").Value() = ""

$wsSteps.Range("B32").Value() = "TC003_13"
$wsSteps.Range("C32").Value() = "Verify quantity of the product"
$wsSteps.Range("D32").Value() = "CartPage"
$wsSteps.Range("E32").Value() = "txtbx_Qty"
$wsSteps.Range("F32").Value() = "verify"
$wsSteps.Range("G32").Value() = ""

$wsSteps.Range("A33").Value() = "TC004"
$wsSteps.Range("B33").Value() = "TC003_14"
$wsSteps.Range("C33").Value() = "Verify price of the product"
$wsSteps.Range("D33").Value() = "CartPage"
$wsSteps.Range("E33").Value() = "lbl_SubTotal"
$wsSteps.Range("F33").Value() = "verify"
$wsSteps.Range("G33").Value() = ""

# Data validation lists on D/E/F need their ranges stretched to row 33 (D's
# sqref already covered through row 19... it is only E and F that grow).
$wsSteps.Range("D2:D19").Validation.Delete()
$wsSteps.Range("E2:E19").Validation.Delete()
$wsSteps.Range("F2:F19").Validation.Delete()
$wsSteps.Range("D2:D19").Validation.Add(3, 1, 1, "PageName")
$wsSteps.Range("E2:E33").Validation.Add(3, 1, 1, "INDIRECT(D2)")
$wsSteps.Range("F2:F33").Validation.Add(3, 1, 1, "ActionKeywords")

# Column C widens (and becomes a "best fit" column) now that it holds the
# longer "Verify and store product price" / "Verify quantity of the
# product" text.
$wsSteps.Columns.Item(3).ColumnWidth = 36.16

# ---------------------------------------------------------------------------
# Sheet "PageObjectModel": add the new "storeValue" action keyword
# ---------------------------------------------------------------------------
$wsPom = $wb.Worksheets.Item("PageObjectModel")
$wsPom.Range("F10").Copy($wsPom.Range("F11"))
$wsPom.Range("F11").Value() = "storeValue"

# ---------------------------------------------------------------------------
# Selections (set last so TestCases ends up the active/tabSelected sheet,
# matching the saved workbook view).
# ---------------------------------------------------------------------------
$wsSteps.Range("H31").Select()
$wsPom.Range("F2").Select()
$wsCases.Range("D4").Select()
